$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells on row 1: "Cotizacion N°" label and its value
$ws.Range("C1").Value = "Cotizacion N°"
$ws.Range("C1").Style = "label_style"
$ws.Range("D1").Value = 1000

# Move the active selection to D6 (no longer A1:B1)
$ws.Range("D6").Select()
